$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Jordan Poole", "PG,SG", "Washington Wizards"),
    @("Kyrie Irving", "PG,SG", "Dallas Mavericks"),
    @("Shai Gilgeous-Alexander", "PG,SG", "Oklahoma City Thunder"),
    @("CJ McCollum", "PG,SG", "New Orleans Pelicans"),
    @("Jalen Williams", "SG,SF,PF,C", "Oklahoma City Thunder"),
    @("Tobias Harris", "SF,PF", "Detroit Pistons"),
    @("Rui Hachimura", "SF,PF", "Los Angeles Lakers"),
    @("Christian Braun", "SG,SF", "Denver Nuggets"),
    @("Jimmy Butler", "SF,PF", "Miami Heat"),
    @("Joel Embiid", "C", "Philadelphia 76ers"),
    @("Guerschon Yabusele", "PF,C", "Philadelphia 76ers"),
    @("Keyonte George", "PG,SG", "Utah Jazz"),
    @("Lauri Markkanen", "SF,PF", "Utah Jazz"),
    @("Zach LaVine", "SG,SF", "Chicago Bulls"),
    @("Jordan Clarkson", "SG,SF", "Utah Jazz"),
    @("John Collins", "PF,C", "Utah Jazz"),
    @("RJ Barrett", "SG,SF,PF", "Toronto Raptors")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
